$wb = $excel.ActiveWorkbook

# Sheet "V": A1 changes from "monetary" to "mixed"
$wsV = $wb.Worksheets.Item("V")
$wsV.Range("A1").Value = "mixed"

# Sheet "U+e": A1 changes from "monetary" to "mixed"; selection resets to A1
$wsUe = $wb.Worksheets.Item("U+e")
$wsUe.Range("A1").Value = "mixed"
$wsUe.Range("A1").Select()

# Sheet "F": A1 stays "physical"; selection resets to A1
$wsF = $wb.Worksheets.Item("F")
$wsF.Range("A1").Value = "physical"
$wsF.Range("A1").Select()

# Keep "V" as the active/selected sheet (unchanged from before the edit)
$wsV.Activate()
